$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Purple"
$ws.Range("B9").Value = "Lion"
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 2
$ws.Range("E9").Formula = "=SUM(C9:D9)"

$ws.Range("E9").Select()
